$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The three data rows (2,3,4) are cyclically rotated:
#   new row2 = old row3
#   new row3 = old row4
#   new row4 = old row2
# Columns A,B,C,E,F,G,H,I,J,K,R stay identical across all three rows,
# so only D,L,M,N,O,P,Q,S,T actually change values.

$oldRow2 = @{
    D = $ws.Range("D2").Value2
    L = $ws.Range("L2").Value2
    M = $ws.Range("M2").Value2
    N = $ws.Range("N2").Value2
    O = $ws.Range("O2").Value2
    P = $ws.Range("P2").Value2
    Q = $ws.Range("Q2").Value2
    S = $ws.Range("S2").Value2
    T = $ws.Range("T2").Value2
}

$oldRow3 = @{
    D = $ws.Range("D3").Value2
    L = $ws.Range("L3").Value2
    M = $ws.Range("M3").Value2
    N = $ws.Range("N3").Value2
    O = $ws.Range("O3").Value2
    P = $ws.Range("P3").Value2
    Q = $ws.Range("Q3").Value2
    S = $ws.Range("S3").Value2
    T = $ws.Range("T3").Value2
}

$oldRow4 = @{
    D = $ws.Range("D4").Value2
    L = $ws.Range("L4").Value2
    M = $ws.Range("M4").Value2
    N = $ws.Range("N4").Value2
    O = $ws.Range("O4").Value2
    P = $ws.Range("P4").Value2
    Q = $ws.Range("Q4").Value2
    S = $ws.Range("S4").Value2
    T = $ws.Range("T4").Value2
}

function Set-RowValues($row, $data) {
    $ws.Range("D$row").Value2 = $data.D
    $ws.Range("L$row").Value2 = $data.L
    $ws.Range("M$row").Value2 = $data.M
    $ws.Range("N$row").Value2 = $data.N
    $ws.Range("O$row").Value2 = $data.O
    $ws.Range("P$row").Value2 = $data.P
    $ws.Range("Q$row").Value2 = $data.Q
    $ws.Range("S$row").Value2 = $data.S
    $ws.Range("T$row").Value2 = $data.T
}

Set-RowValues 2 $oldRow3
Set-RowValues 3 $oldRow4
Set-RowValues 4 $oldRow2
